$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) retains text formatting so numeric-looking
# strings (e.g. "331.85") are not coerced into numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '27.733.21'
$ws.Range('D3').Value = '1.884.77'
$ws.Range('E3').Value = '  -0.17%  '
$ws.Range('E4').Value = '  -0.28%  '
$ws.Range('D5').Value = '331.85'
$ws.Range('E5').Value = '  +1.85%  '
$ws.Range('E6').Value = '  -0.33%  '
$ws.Range('D7').Value = '0.4726'
$ws.Range('E7').Value = '  +3.02%  '
$ws.Range('D8').Value = '0.3978'
$ws.Range('E8').Value = '  +1.19%  '
$ws.Range('D9').Value = '48.32'
$ws.Range('E9').Value = '  -5.91%  '
$ws.Range('D10').Value = '0.08060'
$ws.Range('E10').Value = '  -2.00%  '
$ws.Range('D11').Value = '1.027'
$ws.Range('E11').Value = '  -0.79%  '
$ws.Range('D12').Value = '21.87'
$ws.Range('E12').Value = '  +1.48%  '
$ws.Range('D13').Value = '1.883.32'
$ws.Range('E13').Value = '  -0.64%  '
$ws.Range('D14').Value = '5.972'
$ws.Range('E14').Value = '  +0.06%  '
$ws.Range('D15').Value = '7.208'
$ws.Range('E15').Value = '  -1.31%  '
$ws.Range('D16').Value = '1.001'
$ws.Range('E16').Value = '  -0.41%  '
$ws.Range('D17').Value = '87.02'
$ws.Range('E17').Value = '  -2.26%  '
$ws.Range('D18').Value = '0.00001044'
$ws.Range('E18').Value = '  -1.32%  '
$ws.Range('D19').Value = '0.06598'
$ws.Range('E19').Value = '  +0.32%  '
$ws.Range('D20').Value = '17.25'
$ws.Range('E20').Value = '  -1.72%  '
$ws.Range('D21').Value = '1.001'
$ws.Range('E21').Value = '  -0.40%  '
$ws.Range('D22').Value = '27.750.02'
$ws.Range('E22').Value = '  -0.94%  '
$ws.Range('D23').Value = '5.518'
$ws.Range('E23').Value = '  -2.12%  '
$ws.Range('E24').Value = '  -0.77%  '
$ws.Range('D25').Value = '2.307'
$ws.Range('E25').Value = '  -0.05%  '
$ws.Range('D26').Value = '2.099.34'
$ws.Range('E26').Value = '  -2.63%  '
$ws.Range('D27').Value = '155.08'
$ws.Range('E27').Value = '  +0.65%  '
$ws.Range('D28').Value = '20.22'
$ws.Range('E28').Value = '  +1.72%  '
$ws.Range('D29').Value = '2.103'
$ws.Range('E29').Value = '  +0.15%  '
$ws.Range('D30').Value = '5.600'
$ws.Range('E30').Value = '  -0.97%  '
$ws.Range('D31').Value = '122.61'
$ws.Range('E31').Value = '  -1.03%  '
$ws.Range('D32').Value = '0.9690'
$ws.Range('E32').Value = '  +1.23%  '
$ws.Range('D33').Value = '0.09539'
$ws.Range('E33').Value = '  +0.08%  '
$ws.Range('D34').Value = '1.474'
$ws.Range('E34').Value = '  +0.97%  '
$ws.Range('D35').Value = '3.621'
$ws.Range('E35').Value = '  -0.27%  '
$ws.Range('D36').Value = '5.308'
$ws.Range('E36').Value = '  -2.72%  '
$ws.Range('D37').Value = '0.06129'
$ws.Range('E37').Value = '  +0.47%  '
$ws.Range('D38').Value = '0.02257'
$ws.Range('E38').Value = '  -0.99%  '
$ws.Range('E39').Value = '  -1.87%  '
$ws.Range('D40').Value = '8.184'
$ws.Range('E40').Value = '  -5.31%  '
$ws.Range('D41').Value = '0.6015'
$ws.Range('E41').Value = '  -1.24%  '
$ws.Range('E42').Value = '  -0.32%  '
$ws.Range('D43').Value = '0.1900'
$ws.Range('E43').Value = '  +0.86%  '
$ws.Range('D44').Value = '10.35'
$ws.Range('E44').Value = '  -3.17%  '
$ws.Range('D45').Value = '1.252'
$ws.Range('E45').Value = '  -4.21%  '
$ws.Range('E46').Value = '  -1.84%  '
$ws.Range('D47').Value = '12.25'
$ws.Range('E47').Value = '  -3.09%  '
$ws.Range('D48').Value = '3.406'
$ws.Range('E48').Value = '  -0.52%  '
$ws.Range('D49').Value = '1.942'
$ws.Range('E49').Value = '  -2.37%  '
$ws.Range('D50').Value = '0.06833'
$ws.Range('E50').Value = '  -0.67%  '
$ws.Range('D51').Value = '110.65'
$ws.Range('E51').Value = '  +0.43%  '
